$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.302.29"
$ws.Range("D3").Value = "1.929.88"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.66"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7201"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.93"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3203"
$ws.Range("E9").Value = "  -4.26%  "
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7883"
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08016"
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "1.931.22"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.70"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.65"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "30.296.80"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "256.78"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008058"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.731"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").Value = "2.183.06"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9982"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.555"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.48"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.300"
$ws.Range("E28").Value = "  -4.78%  "
$ws.Range("E29").Value = "  -2.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.353"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.533"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.420"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.152"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05111"
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.289"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7495"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01985"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.799"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.45"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.402"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4522"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.996"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8456"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.38"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.841"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.486"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.85"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "962.30"
$ws.Range("E50").Value = "  +8.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4210"
$ws.Range("E51").Value = "  +0.21%  "
